# Update the "Updating personal info user and operator" task row on the
# "AD Gantt Chart & Burndown" sheet: fix the task title text and correct the
# estimate/completed hours (task is actually fully completed: 6/6 not 5/8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AD Gantt Chart & Burndown")

# Fix task title wording: "... personal info user and operator" -> "... personal info of user and operator"
$ws.Range("C45").Value = "Updating personal info of user and operator"

# Correct hours: estimate 8 -> 6, completed 5 -> 6 (task is fully done)
$ws.Range("E45").Value = 6
$ws.Range("F45").Value = 6

# Move the active selection near the edited cells, as left by the editor.
$ws.Activate()
$ws.Range("H46").Select()
